# Update "Tasas on-shore 2021 - Diaria" with the latest daily observations
# (MV data refresh): four new trading days appended after 08-09-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append at the bottom of the table.
$newRows = @(
    @{ Row = 175; Fecha = "09-09-2021"; B = -0.19; C = -0.01; D = 0.04 },
    @{ Row = 176; Fecha = "10-09-2021"; B = -0.13; C = 0.04;  D = 0.23 },
    @{ Row = 177; Fecha = "13-09-2021"; B = 0.02;  C = 0.14;  D = 0.42 },
    @{ Row = 178; Fecha = "14-09-2021"; B = 0.02;  C = 0.27;  D = 0.35 }
)

# Force column A to be treated as plain text so the dd-mm-yyyy labels are
# not auto-converted into Excel date serials (matching the rest of the
# "Serie" column, which stores these as text).
$labelRange = $ws.Range("A175:A178")
$labelRange.NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Fecha
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}

# Drop the direct text-format override so the new cells keep the sheet's
# default (General) style, same as the pre-existing "Serie" entries.
$labelRange.ClearFormats()
